# Reset the payroll template back to a blank state so it can be reused for
# the next pay period: the three "expected date" placeholders are emptied,
# the now-unused helper cells around them are wiped completely, and the
# sample grand-total row at the bottom is cleared out too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expected Period Ending / Pay / Run dates (K1:K3) - blank the value but
# keep the date number-format applied to the cells.
$ws.Range("K1:K3").ClearContents()

# I1:J3 (blank spacer cells) and the whole helper row 4 (I4:K4) are no
# longer needed once the dates are gone - clear them fully so they drop
# out of the sheet entirely.
$ws.Range("I1:J3").Clear()
$ws.Range("I4:K4").Clear()

# Clear the sample "TOTALS" values left over from testing (B25:H25); the
# remaining totals cells (I25:K25) were already blank.
$ws.Range("B25:H25").ClearContents()

# Leave the cursor on K3 (the period-ending date cell) ready for entry.
$ws.Range("K3").Select()

$wb.Save()
